$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 10 (pushing existing rows 10+ down to 12+).
# This naturally carries down the existing cell styles (including the
# alternating "banding" fill on column C), matching the target layout.
$ws.Range("A10:A11").EntireRow.Insert()

# Fill in the two new requirement rows.
$ws.Range("A10").Value = "09FU"
$ws.Range("B10").Value = "Deve ser permitido que o usuário remova uma oferta de doação"
$ws.Range("C10").Value = "Funcional"

$ws.Range("A11").Value = "10FU"
$ws.Range("B11").Value = "Deve ser permitido que a instituição remova um pedido de doação"
$ws.Range("C11").Value = "Funcional"

# Renumber the requirement codes for all the rows that shifted down,
# keeping the sequence contiguous (11FU..18FU, then 19NF..22NF).
$ws.Range("A12").Value = "11FU"
$ws.Range("A13").Value = "12FU"
$ws.Range("A14").Value = "13FU"
$ws.Range("A15").Value = "14FU"
$ws.Range("A16").Value = "15FU"
$ws.Range("A17").Value = "16FU"
$ws.Range("A18").Value = "17FU"
$ws.Range("A19").Value = "18FU"
$ws.Range("A20").Value = "19NF"
$ws.Range("A21").Value = "20NF"
$ws.Range("A22").Value = "21NF"
$ws.Range("A23").Value = "22NF"
